# Update TPM-derived values in Sheet1 per new script run.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 0.5431923333333334
$ws.Range("H2").Value = 1.629577
$ws.Range("M2").Value = 9.278280333333333
$ws.Range("N2").Value = 27.834841
$ws.Range("O2").Value = 0.2122966588143784
$ws.Range("P2").Value = 0.2122966588143784
$ws.Range("Q2").Value = 5.039890743584111
$ws.Range("R2").Value = 45.359016692257
$ws.Range("S2").Value = 0.2122966588143784
$ws.Range("T2").Value = 0.2122966588143784

# Row 3
$ws.Range("G3").Value = 0.5431923333333334
$ws.Range("H3").Value = 1.629577
$ws.Range("O3").Value = 0.2154323368929792
$ws.Range("P3").Value = 0.2154323368929792
$ws.Range("Q3").Value = 5.114331269457
$ws.Range("R3").Value = 46.028981425113
$ws.Range("S3").Value = 0.2154323368929792
$ws.Range("T3").Value = 0.2154323368929792

# Row 4
$ws.Range("G4").Value = 0.5431923333333334
$ws.Range("H4").Value = 1.629577
$ws.Range("M4").Value = 7.033255
$ws.Range("N4").Value = 21.099765
$ws.Range("O4").Value = 0.1609281551588013
$ws.Range("P4").Value = 0.1609281551588013
$ws.Range("Q4").Value = 3.820410194378334
$ws.Range("R4").Value = 34.383691749405
$ws.Range("S4").Value = 0.1609281551588013
$ws.Range("T4").Value = 0.1609281551588013

# Row 5
$ws.Range("G5").Value = 0.5431923333333334
$ws.Range("H5").Value = 1.629577
$ws.Range("M5").Value = 17.977458
$ws.Range("N5").Value = 53.932374
$ws.Range("O5").Value = 0.4113428491338411
$ws.Range("P5").Value = 0.411342849133841
$ws.Range("Q5").Value = 9.765217358421999
$ws.Range("R5").Value = 87.886956225798
$ws.Range("S5").Value = 0.4113428491338411
$ws.Range("T5").Value = 0.411342849133841
